$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new daily rows (四方坪站 / 高岭站) for 2025-11-25 (serial 45986)
$ws.Range("A50").Value = 45986
$ws.Range("B50").Value = "四方坪站"
$ws.Range("C50").Value = 8144.57
$ws.Range("D50").Value = 7219.49
$ws.Range("E50").Value = 2728.64
$ws.Range("F50").Value = 357

$ws.Range("A51").Value = 45986
$ws.Range("B51").Value = "高岭站"
$ws.Range("C51").Value = 4037.14
$ws.Range("D51").Value = 3554.79
$ws.Range("E51").Value = 999.4
$ws.Range("F51").Value = 156

# Match the author's view state: selection moved to L49
$ws.Range("L49").Select()
